# Apply updated odds values for rows 7, 8, 11, and 12 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G7").Value = 2.12
$ws.Range("H7").Value = 2.87
$ws.Range("I7").Value = 3.55
$ws.Range("J7").Value = 2.77
$ws.Range("K7").Value = 1.95
$ws.Range("L7").Value = 4.3
$ws.Range("S7").Value = 1.52
$ws.Range("T7").Value = 2.42
$ws.Range("W7").Value = 5.9
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 21
$ws.Range("AA7").Value = 20
$ws.Range("AB7").Value = 37
$ws.Range("AD7").Value = 5.8
$ws.Range("AE7").Value = 17.5
$ws.Range("AG7").Value = 7.9
$ws.Range("AH7").Value = 17.5
$ws.Range("AI7").Value = 13
$ws.Range("AJ7").Value = 55
$ws.Range("AK7").Value = 45
$ws.Range("AL7").Value = 60
$ws.Range("AN7").Value = 3.85
$ws.Range("AO7").Value = 11.25
$ws.Range("AP7").Value = 22
$ws.Range("AQ7").Value = 45
$ws.Range("AR7").Value = 90
$ws.Range("AS7").Value = 2.37
$ws.Range("AT7").Value = 7.5
$ws.Range("AV7").Value = 5.4
$ws.Range("AW7").Value = 22
$ws.Range("AX7").Value = 32
$ws.Range("AZ7").Value = 200
$ws.Range("BA7").Value = 500
$ws.Range("BB7").Value = 350
$ws.Range("G8").Value = 1.57
$ws.Range("H8").Value = 3.75
$ws.Range("I8").Value = 6.25
$ws.Range("L8").Value = 6
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("Q8").Value = 1.95
$ws.Range("R8").Value = 1.9
$ws.Range("X8").Value = 7
$ws.Range("Z8").Value = 11
$ws.Range("AJ8").Value = 67
$ws.Range("AL8").Value = 51
$ws.Range("AN8").Value = 3.5
$ws.Range("AO8").Value = 8
$ws.Range("BA8").Value = 301
$ws.Range("G11").Value = 1.93
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 3.9
$ws.Range("J11").Value = 2.52
$ws.Range("K11").Value = 2.02
$ws.Range("L11").Value = 4.3
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 8
$ws.Range("O11").Value = 1.36
$ws.Range("P11").Value = 2.67
$ws.Range("Q11").Value = 2.05
$ws.Range("R11").Value = 1.62
$ws.Range("S11").Value = 1.42
$ws.Range("T11").Value = 2.45
$ws.Range("U11").Value = 1.85
$ws.Range("V11").Value = 1.75
$ws.Range("W11").Value = 6.3
$ws.Range("X11").Value = 8.5
$ws.Range("Y11").Value = 8.5
$ws.Range("Z11").Value = 16.5
$ws.Range("AA11").Value = 17
$ws.Range("AB11").Value = 32
$ws.Range("AC11").Value = 8
$ws.Range("AD11").Value = 6.2
$ws.Range("AE11").Value = 16
$ws.Range("AF11").Value = 80
$ws.Range("AG11").Value = 10
$ws.Range("AH11").Value = 21
$ws.Range("AI11").Value = 13
$ws.Range("AJ11").Value = 60
$ws.Range("AK11").Value = 40
$ws.Range("AL11").Value = 50
$ws.Range("AM11").Value = 800
$ws.Range("AN11").Value = 3.7
$ws.Range("AO11").Value = 9.75
$ws.Range("AP11").Value = 19.5
$ws.Range("AQ11").Value = 37
$ws.Range("AR11").Value = 75
$ws.Range("AS11").Value = 2.42
$ws.Range("AT11").Value = 7.3
$ws.Range("AU11").Value = 70
$ws.Range("AV11").Value = 5.6
$ws.Range("AW11").Value = 22
$ws.Range("AX11").Value = 29
$ws.Range("AY11").Value = 120
$ws.Range("AZ11").Value = 175
$ws.Range("BA11").Value = 400
$ws.Range("G12").Value = 2.7
$ws.Range("I12").Value = 2.52
$ws.Range("J12").Value = 3.2
$ws.Range("L12").Value = 3.1
$ws.Range("N12").Value = 10
$ws.Range("W12").Value = 8.5
$ws.Range("Y12").Value = 10
$ws.Range("Z12").Value = 32
$ws.Range("AA12").Value = 23
$ws.Range("AB12").Value = 32
$ws.Range("AG12").Value = 7.8
$ws.Range("AI12").Value = 9.5
$ws.Range("AJ12").Value = 28
$ws.Range("AK12").Value = 22
$ws.Range("AL12").Value = 32
$ws.Range("AN12").Value = 4.6
$ws.Range("AO12").Value = 14
$ws.Range("AP12").Value = 20
$ws.Range("AQ12").Value = 60
$ws.Range("AR12").Value = 90
$ws.Range("AT12").Value = 6.7
$ws.Range("AV12").Value = 4.4
$ws.Range("AW12").Value = 13.5
$ws.Range("AX12").Value = 21
$ws.Range("AY12").Value = 55
$ws.Range("AZ12").Value = 90
